# preparation publication 0.2.0 be6a807bbdadc24333e2c553161780cb6e805524
#
# Updates the "Metadata" sheet of the StructureDefinition-eclaire-label
# workbook for the 0.2.0 publication:
#   - Version bumped from 0.1.1 to 0.2.0
#   - Date updated to the new publication timestamp
#   - A new "Jurisdiction" / "iso:code:3166:FR" row inserted right after
#     the "Contact" row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the published version number.
$ws.Range("B3").Value = "0.2.0"

# Refresh the publication date/time.
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), pushing
# Description/Purpose/Copyright/etc. down by one row.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# The freshly inserted row picks up a default "no border" style from the
# Insert() call; copy the formatting from the row below (still using the
# standard body style) so the new row matches the rest of the table.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
